$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Espana and Argentina order (Espana now precedes Argentina in the country list)
$ws.Range("A8").Value = "España"
$ws.Range("A9").Value = "Argentina"

# Update "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 19:32"

# Update per-country statistics
$ws.Range("B4").Value = 8398560
$ws.Range("C4").Value = 10761
$ws.Range("D4").Value = 5465288
$ws.Range("E4").Value = 2708388
$ws.Range("G4").Value = 154
$ws.Range("H4").Value = 224884

$ws.Range("B5").Value = 7590514
$ws.Range("C5").Value = 42276
$ws.Range("D5").Value = 6725219
$ws.Range("E5").Value = 750132
$ws.Range("G5").Value = 521
$ws.Range("H5").Value = 115163

$ws.Range("B6").Value = 5237961
$ws.Range("C6").Value = 2617
$ws.Range("E6").Value = 433949
$ws.Range("G6").Value = 77
$ws.Range("H6").Value = 153982

$ws.Range("B8").Value = 1015795
$ws.Range("C8").Value = 12214
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("G8").Value = 73
$ws.Range("H8").Value = 33992

$ws.Range("B9").Value = 989680
$ws.Range("D9").Value = 803965
$ws.Range("E9").Value = 159448
$ws.Range("H9").Value = 26267

$ws.Range("B21").Value = 371666
$ws.Range("C21").Value = 4685
$ws.Range("E21").Value = 69881
$ws.Range("G21").Value = 19
$ws.Range("H21").Value = 9885

$ws.Range("B24").Value = 349519
$ws.Range("C24").Value = 2026
$ws.Range("D24").Value = 305427
$ws.Range("E24").Value = 34721
$ws.Range("G24").Value = 75
$ws.Range("H24").Value = 9371

$ws.Range("B27").Value = 304635
$ws.Range("C27").Value = 1526
$ws.Range("D27").Value = 272768
$ws.Range("E27").Value = 29607
$ws.Range("G27").Value = 51
$ws.Range("H27").Value = 2260

$ws.Range("B31").Value = 199890
$ws.Range("C31").Value = 1742
$ws.Range("D31").Value = 168689
$ws.Range("E31").Value = 21431
$ws.Range("G31").Value = 10
$ws.Range("H31").Value = 9770

$ws.Range("B35").Value = 175749
$ws.Range("C35").Value = 2117
$ws.Range("D35").Value = 146421
$ws.Range("E35").Value = 26352
$ws.Range("G35").Value = 48
$ws.Range("H35").Value = 2976

$ws.Range("B63").Value = 62944
$ws.Range("C63").Value = 658
$ws.Range("D63").Value = 28855
$ws.Range("E63").Value = 33563
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 526

$ws.Range("B69").Value = 50993
$ws.Range("C69").Value = 1031
$ws.Range("E69").Value = 25777

$ws.Range("B104").Value = 12326
$ws.Range("C104").Value = 33
$ws.Range("D104").Value = 10426
$ws.Range("E104").Value = 1769

$ws.Range("B189").Value = 268
$ws.Range("C189").Value = 3
$ws.Range("D189").Value = 222
$ws.Range("E189").Value = 44

